# Revert "feat(dialog): update CN data and dialogue Excel files"
# This reverts the addition of two "modAffinity" adjustment rows (27 and 35)
# in the big_sister.xlsx dialogue sheet, shifting subsequent rows back up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 35 first (the modAffinity / -200 row), then row 27
# (the modAffinity / 100 row), deleting from the bottom up so row
# numbers for the earlier deletion remain valid.
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(27).Delete()

# Row deletion shifts the sheet-bottom-anchored MAX() range down;
# restore it to its original extent (H1048576 is the literal last row).
$ws.Range("H2").Formula = "=MAX(H4:H1048576)"

# Restore the view state (pane/selection) to match the reverted layout:
# top pane selection K33, frozen at row 2 with the bottom pane scrolled to
# A3, and the bottom pane's active selection at J14.
$ws.Range("K33").Select()
$ws.Range("A3").Select()
$ws.Range("J14").Select()
